$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: becomes a fully-populated "real" entry (issue #7), matching the
#     green-filled look of rows 4/6/7/8. Copy formatting from row 8 first,
#     then overwrite the values/text that differ.
$ws.Range("A8:J8").Copy()
[void]$ws.Range("A10:J10").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B10").Value = 1
$ws.Range("C10").Value = "FIX"
$ws.Range("D10").Value = "1.3.3"
$ws.Range("E10").Value = "1.3.3"
$ws.Range("F10").Value = "STRIPS"
$ws.Range("G10").Value = "CONT mode"
$ws.Range("I10").Value = "FIX IT!"
$ws.Range("J10").Value = "Reverted code changes that caused the issue.  I am an idiot"
$ws.Range("H10").Value = "Timing runnning slow per strip"
$ws.Rows(10).RowHeight = 30

# --- Row 11: stays unfilled/white, only a handful of cells get values.
$ws.Range("B11").Value = 1
$ws.Range("D11").Value = "1.3.3"
$ws.Range("F11").Value = "STRIPS"
$ws.Range("G11").Value = "CONT mode"
$ws.Range("H11").Value = "FIxing 7 above broke [RST] in STRIPS mode"
$ws.Range("I11").Value = "Restore functionality"

# --- Clear the stale selection marker left over from the original file so
#     the saved sheet view doesn't keep pointing at C10.
[void]$ws.Range("A1").Select()

Write-Output "edit complete"
